$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "69.057.67", "0.0000161").
# Force a Text number format on just the cells we touch so Excel does not
# reinterpret them as numbers (which would mangle values like "0.0000161"
# into scientific notation, or "1.00" into 1).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values cell by cell, matching the authored diff.
$ws.Range("D2").Value = "69.057.67"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "3.808.96"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").Value = "600.91"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "163.87"
$ws.Range("E6").Value = "  -2.48%  "
$ws.Range("D7").Value = "3.809.06"
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("E10").Value = "  +1.77%  "
$ws.Range("D11").Value = "6.31"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "37.12"
$ws.Range("E13").Value = "  -2.40%  "
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "4.439.28"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").Value = "3.815.35"
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("D17").Value = "69.150.75"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "7.47"
$ws.Range("E18").Value = "  +2.38%  "
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "11.49"
$ws.Range("E20").Value = "  +5.61%  "
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("D22").Value = "485.92"
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("D23").Value = "0.720"
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("D24").Value = "0.0000161"
$ws.Range("E24").Value = "  +7.11%  "
$ws.Range("D25").Value = "84.58"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("E26").Value = "  -2.37%  "
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "2.97"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("D31").Value = "8.02"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("E32").Value = "  -4.71%  "
$ws.Range("D33").Value = "3.961.20"
$ws.Range("E33").Value = "  +1.90%  "
$ws.Range("D34").Value = "31.86"
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("D35").Value = "3.749.04"
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("E36").Value = "  -1.57%  "
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").Value = "  +4.63%  "
$ws.Range("D39").Value = "5.87"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "3.04"
$ws.Range("E41").Value = "  +1.59%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Value = "0.319"
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("D43").Value = "437.69"
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("E47").Value = "  -1.18%  "
$ws.Range("D48").Value = "2.828.06"
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("D49").Value = "142.14"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").Value = "39.28"
$ws.Range("E50").Value = "  -2.58%  "
$ws.Range("D51").Value = "0.0352"
$ws.Range("E51").Value = "  -0.49%  "
